$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "2018-11-01 00:00:00"
$ws.Range("C2").Value = 2785

$ws.Range("B3").Value = "2018-11-02 00:00:00"

$ws.Range("B4").Value = "2018-11-03 00:00:00"

$ws.Range("B5").Value = "2018-11-04 00:00:00"
$ws.Range("C5").Value = 100

$ws.Range("B6").Value = "2018-11-05 00:00:00"
$ws.Range("C6").Value = 5
